$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: advance the weekly report by one week ---
# "Volume 31   Number  50" -> "Volume 31   Number  51"
$ws.Range("A8").Value = "Volume 31   Number  51"
# "Report Covering the Week  12/9/2024  Through  12/15/2024" -> next week's range
$ws.Range("C9").Value = "Report Covering the Week  12/16/2024  Through  12/22/2024"

# --- Insert a new blank row before the footer block (old row 56 -> 57, old row 57 -> 58) ---
$ws.Rows(56).Insert()
$ws.Range("A56").Clear()

# --- Bulk numeric data refresh across the crime-complaint table (rows 14-31, 33) ---
$data = @{
    "F14" = 4
    "H14" = -33.333333333333
    "I14" = 49
    "K14" = -18.333333333333
    "L14" = -19.672131147541
    "M14" = -43.678160919540
    "N14" = -79.668049792531
    "C15" = 2
    "D15" = 3
    "E15" = -33.333333333333
    "F15" = 10
    "G15" = 20
    "H15" = -50
    "I15" = 219
    "J15" = 219
    "K15" = 0
    "L15" = 4.285714285714
    "M15" = 27.325581395348
    "N15" = -59.369202226345
    "C16" = 24
    "D16" = 37
    "E16" = -35.135135135135
    "F16" = 108
    "G16" = 131
    "H16" = -17.557251908396
    "I16" = 1659
    "J16" = 1737
    "K16" = -4.490500863557
    "L16" = -13.322884012539
    "M16" = -41.687170474516
    "N16" = -87.902872976520
    "C17" = 51
    "D17" = 69
    "E17" = -26.086956521739
    "F17" = 234
    "G17" = 272
    "H17" = -13.970588235294
    "I17" = 3585
    "J17" = 3461
    "K17" = 3.582779543484
    "L17" = 7.367475292003
    "M17" = 48.447204968944
    "N17" = -45.158329508949
    "C18" = 38
    "D18" = 38
    "E18" = 0
    "F18" = 105
    "G18" = 115
    "H18" = -8.695652173913
    "I18" = 1486
    "J18" = 1655
    "K18" = -10.211480362537
    "L18" = -29.573459715639
    "M18" = -56.100443131462
    "N18" = -91.901907356948
    "C19" = 77
    "D19" = 122
    "E19" = -36.885245901639
    "F19" = 350
    "G19" = 488
    "H19" = -28.278688524590
    "I19" = 5548
    "J19" = 6382
    "K19" = -13.068003760576
    "L19" = -21.081081081081
    "M19" = 1.630335226231
    "N19" = -36.731668377238
    "C20" = 40
    "D20" = 34
    "E20" = 17.647058823529
    "F20" = 144
    "H20" = -4
    "I20" = 1997
    "J20" = 1833
    "K20" = 8.947081287506
    "L20" = 13.659647125782
    "M20" = 6.905781584582
    "N20" = -91.211934518570
    "C21" = 234
    "D21" = 303
    "E21" = -22.772277227722
    "F21" = 955
    "G21" = 1182
    "H21" = -19.204737732656
    "I21" = 14543
    "J21" = 15347
    "K21" = -5.238808887730
    "L21" = -11.436575117227
    "M21" = -10.399852134803
    "N21" = -79.480486497164
    "C22" = 6
    "D22" = 2
    "E22" = 200
    "F22" = 20
    "G22" = 9
    "H22" = 122.222222222222
    "I22" = 199
    "J22" = 198
    "K22" = 0.505050505050
    "L22" = 1.530612244897
    "M22" = -27.372262773722
    "C23" = 8
    "D23" = 10
    "E23" = -20
    "G23" = 46
    "H23" = -19.565217391304
    "I23" = 488
    "J23" = 517
    "K23" = -5.609284332688
    "L23" = -11.913357400722
    "M23" = 45.238095238095
    "C24" = 286
    "E24" = -3.050847457627
    "F24" = 1154
    "G24" = 1137
    "H24" = 1.495162708883
    "I24" = 14691
    "J24" = 15386
    "K24" = -4.517093461588
    "L24" = -8.683490800596
    "M24" = 21.413223140495
    "D25" = 134
    "E25" = -15.671641791044
    "F25" = 489
    "G25" = 520
    "H25" = -5.961538461538
    "I25" = 6923
    "J25" = 7411
    "K25" = -6.584806368911
    "L25" = -9.868506704856
    "C26" = 108
    "D26" = 128
    "E26" = -15.625
    "F26" = 477
    "G26" = 483
    "H26" = -1.242236024844
    "I26" = 6395
    "J26" = 5787
    "K26" = 10.506307240366
    "L26" = 18.338267949666
    "M26" = -1.856967464702
    "C27" = 3
    "D27" = 5
    "E27" = -40
    "F27" = 15
    "G27" = 29
    "H27" = -48.275862068965
    "I27" = 314
    "J27" = 320
    "K27" = -1.875
    "L27" = -6.824925816023
    "C28" = 13
    "D28" = 14
    "E28" = -7.142857142857
    "G28" = 34
    "H28" = 0
    "I28" = 659
    "J28" = 620
    "K28" = 6.290322580645
    "L28" = -3.654970760233
    "C29" = 2
    "D29" = 5
    "E29" = -60
    "F29" = 9
    "G29" = 19
    "H29" = -52.631578947368
    "I29" = 106
    "J29" = 153
    "K29" = -30.718954248366
    "L29" = -47.783251231527
    "M29" = -61.029411764705
    "N29" = -86.034255599473
    "D30" = 4
    "E30" = -75
    "I30" = 90
    "J30" = 130
    "K30" = -30.769230769230
    "L30" = -43.75
    "M30" = -60.176991150442
    "N30" = -86.425339366515
    "D31" = 4
    "F31" = 3
    "G31" = 16
    "H31" = -81.25
    "I31" = 145
    "J31" = 116
    "K31" = 25
    "L31" = 22.881355932203
    "D33" = 1
    "E33" = -100
    "G33" = 2
    "H33" = 50
    "J33" = 46
    "K33" = -6.521739130434
    "L33" = -6.521739130434
}
foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}

# --- D33/E33 switch from the "n/a" placeholder text to real computed numbers ---
$ws.Range("D33").NumberFormat = "#,##0"
$ws.Range("E33").NumberFormat = "#,##0.0;""-""#,##0.0"
